$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (originally rows 56 and 57). Delete from the
# bottom up so earlier row indices are not shifted while deleting, which
# shrinks the used range from A1:E57 down to A1:E55.
$ws.Rows(57).Delete()
$ws.Rows(56).Delete()

# New B:E values (lamda_1, lamda_2, dic_nbre_clients_poisson_2_keys,
# dic_nbre_clients_prob_poisson_2_values) for data rows 2-55. Column A
# (auto scs index 0..53) is unchanged.
$data = @(
    @(33.94444444444444, 1.95, 0, 0.136),
    @(33.94444444444444, 1.95, 2, 0.001),
    @(33.94444444444444, 1.95, 3, 0.006),
    @(33.94444444444444, 1.95, 4, 0.009000000000000001),
    @(33.94444444444444, 1.95, 5, 0.019),
    @(33.94444444444444, 1.95, 6, 0.038),
    @(33.94444444444444, 1.95, 7, 0.058),
    @(33.94444444444444, 1.95, 8, 0.058),
    @(33.94444444444444, 1.95, 9, 0.04),
    @(33.94444444444444, 1.95, 10, 0.038),
    @(33.94444444444444, 1.95, 11, 0.023),
    @(33.94444444444444, 1.95, 12, 0.02),
    @(33.94444444444444, 1.95, 13, 0.028),
    @(33.94444444444444, 1.95, 14, 0.036),
    @(33.94444444444444, 1.95, 15, 0.023),
    @(33.94444444444444, 1.95, 16, 0.042),
    @(33.94444444444444, 1.95, 17, 0.038),
    @(33.94444444444444, 1.95, 18, 0.03),
    @(33.94444444444444, 1.95, 19, 0.037),
    @(33.94444444444444, 1.95, 20, 0.018),
    @(33.94444444444444, 1.95, 21, 0.02),
    @(33.94444444444444, 1.95, 22, 0.036),
    @(33.94444444444444, 1.95, 23, 0.021),
    @(33.94444444444444, 1.95, 24, 0.018),
    @(33.94444444444444, 1.95, 25, 0.026),
    @(33.94444444444444, 1.95, 26, 0.017),
    @(33.94444444444444, 1.95, 27, 0.014),
    @(33.94444444444444, 1.95, 28, 0.014),
    @(33.94444444444444, 1.95, 29, 0.012),
    @(33.94444444444444, 1.95, 30, 0.016),
    @(33.94444444444444, 1.95, 31, 0.019),
    @(33.94444444444444, 1.95, 32, 0.008),
    @(33.94444444444444, 1.95, 33, 0.009000000000000001),
    @(33.94444444444444, 1.95, 34, 0.007),
    @(33.94444444444444, 1.95, 35, 0.008),
    @(33.94444444444444, 1.95, 36, 0.01),
    @(33.94444444444444, 1.95, 37, 0.004),
    @(33.94444444444444, 1.95, 38, 0.004),
    @(33.94444444444444, 1.95, 39, 0.003),
    @(33.94444444444444, 1.95, 40, 0.004),
    @(33.94444444444444, 1.95, 41, 0.007),
    @(33.94444444444444, 1.95, 42, 0.004),
    @(33.94444444444444, 1.95, 43, 0.001),
    @(33.94444444444444, 1.95, 45, 0.002),
    @(33.94444444444444, 1.95, 47, 0.004),
    @(33.94444444444444, 1.95, 48, 0.002),
    @(33.94444444444444, 1.95, 49, 0.001),
    @(33.94444444444444, 1.95, 50, 0.001),
    @(33.94444444444444, 1.95, 51, 0.002),
    @(33.94444444444444, 1.95, 53, 0.002),
    @(33.94444444444444, 1.95, 58, 0.001),
    @(33.94444444444444, 1.95, 60, 0.002),
    @(33.94444444444444, 1.95, 64, 0.001),
    @(33.94444444444444, 1.95, 65, 0.001)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 2).Value = $row[$j]
    }
}
